$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
# Row 51
$ws.Range("H51").Value = 3958.3333
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3958.3333
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3958.3333
$ws.Range("M51").Value = -4926.3333
$ws.Range("N51").Value = -4926.3333
# Row 69
$ws.Range("H69").Value = 7198.5
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 7198.5
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 113
$ws.Range("H113").Value = 3100
$ws.Range("I113").Value = 2700
$ws.Range("J113").Value = 3900
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 3900
$ws.Range("M113").Value = 554
$ws.Range("N113").Value = -10408
# Row 137
$ws.Range("H137").Value = 2425.4167
$ws.Range("I137").Value = 1648.4706
$ws.Range("K137").Value = 4945.4118
$ws.Range("M137").Value = -2395.4118
# Row 141
$ws.Range("H141").Value = 4743.125
$ws.Range("I141").Value = 3723.3333
$ws.Range("K141").Value = 11169.9999
$ws.Range("M141").Value = -5989.999899999999

$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 4999.5
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 4999
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 4999
$ws.Range("M16").Value = -4713
$ws.Range("N16").Value = -5573
# Row 44
$ws.Range("H44").Value = 68998.336
$ws.Range("J44").Value = 68998.336
$ws.Range("L44").Value = 68998.336
$ws.Range("N44").Value = -69974.336
# Row 45
$ws.Range("H45").Value = 1912.6
$ws.Range("I45").Value = 1897.25
$ws.Range("K45").Value = 1897.25
$ws.Range("M45").Value = -1520.25
# Row 55
$ws.Range("H55").Value = 39992
$ws.Range("I55").Value = 40000
$ws.Range("K55").Value = 40000
$ws.Range("M55").Value = -39685
# Row 61
$ws.Range("H61").Value = 2136.0625
$ws.Range("I61").Value = 2055.6428
$ws.Range("K61").Value = 2055.6428
$ws.Range("M61").Value = -1843.6428
# Row 80
$ws.Range("H80").Value = 89999
$ws.Range("J80").Value = 89999
$ws.Range("L80").Value = 89999
$ws.Range("N80").Value = -91995
# Row 83
$ws.Range("H83").Value = 89999
$ws.Range("J83").Value = 89999
$ws.Range("L83").Value = 269997
$ws.Range("N83").Value = -279981
# Row 97
$ws.Range("H97").Value = 3338.5715
$ws.Range("I97").Value = 1842.5
$ws.Range("J97").Value = 5333.3335
$ws.Range("K97").Value = 1842.5
$ws.Range("L97").Value = 5333.3335
$ws.Range("M97").Value = -1346.5
$ws.Range("N97").Value = -6325.3335
# Row 136
$ws.Range("H136").Value = 2136.0625
$ws.Range("I136").Value = 2055.6428
$ws.Range("K136").Value = 6166.928400000001
$ws.Range("M136").Value = -3616.928400000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 8299.833000000001
$ws.Range("I20").Value = 8299.833000000001
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 8299.833000000001
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -8052.833000000001
$ws.Range("N20").ClearContents()
# Row 107
$ws.Range("H107").Value = 1263.6364
$ws.Range("I107").Value = 1270
$ws.Range("K107").Value = 1270
$ws.Range("M107").Value = 650

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3329.5833
$ws.Range("I58").Value = 1332.5714
$ws.Range("K58").Value = 1332.5714
$ws.Range("M58").Value = -1129.5714
# Row 62
$ws.Range("H62").Value = 48130.332
$ws.Range("I62").Value = 4146.75
$ws.Range("K62").Value = 4146.75
$ws.Range("M62").Value = -3522.75
# Row 65
$ws.Range("H65").Value = 48130.332
$ws.Range("I65").Value = 4146.75
$ws.Range("K65").Value = 20733.75
$ws.Range("M65").Value = -17613.75
# Row 105
$ws.Range("H105").Value = 1807
$ws.Range("I105").Value = 1401
$ws.Range("J105").Value = 2010
$ws.Range("K105").Value = 1401
$ws.Range("L105").Value = 2010
$ws.Range("M105").Value = 346
$ws.Range("N105").Value = -5504
# Row 107
$ws.Range("H107").Value = 742.2
$ws.Range("I107").Value = 703.6667
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 703.6667
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1216.3333
$ws.Range("N107").Value = -4640
# Row 132
$ws.Range("H132").Value = 2067.9429
$ws.Range("I132").Value = 1834.6765
$ws.Range("K132").Value = 5504.029500000001
$ws.Range("M132").Value = -2974.029500000001
# Row 134
$ws.Range("H134").Value = 2230.1904
$ws.Range("I134").Value = 1232.5
$ws.Range("K134").Value = 3697.5
$ws.Range("M134").Value = -1162.5
# Row 136
$ws.Range("H136").Value = 3329.5833
$ws.Range("I136").Value = 1332.5714
$ws.Range("K136").Value = 3997.7142
$ws.Range("M136").Value = -1447.7142

$ws = $wb.Worksheets.Item("CUL")
# Row 140
$ws.Range("H140").Value = 1680.4166
$ws.Range("I140").Value = 1680.4166
$ws.Range("K140").Value = 5041.2498
$ws.Range("M140").Value = 138.7502000000004
# Row 141
$ws.Range("H141").Value = 8007.25
$ws.Range("I141").Value = 8007.25
$ws.Range("K141").Value = 24021.75
$ws.Range("M141").Value = -18841.75

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 2316.6667
$ws.Range("I43").Value = 1480
$ws.Range("J43").Value = 6500
$ws.Range("K43").Value = 1480
$ws.Range("L43").Value = 6500
$ws.Range("M43").Value = -1329
$ws.Range("N43").Value = -6802
# Row 80
$ws.Range("H80").Value = 6318.25
$ws.Range("I80").Value = 5609.6
$ws.Range("J80").Value = 7499.3335
$ws.Range("K80").Value = 5609.6
$ws.Range("L80").Value = 7499.3335
$ws.Range("M80").Value = -4611.6
$ws.Range("N80").Value = -9495.333500000001
# Row 83
$ws.Range("H83").Value = 6318.25
$ws.Range("I83").Value = 5609.6
$ws.Range("J83").Value = 7499.3335
$ws.Range("K83").Value = 28048
$ws.Range("L83").Value = 37496.6675
$ws.Range("M83").Value = -23056
$ws.Range("N83").Value = -47480.6675
# Row 132
$ws.Range("H132").Value = 2763.0908
$ws.Range("I132").Value = 1921.25
$ws.Range("K132").Value = 5763.75
$ws.Range("M132").Value = -3233.75

$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 50021500
$ws.Range("I13").Value = 100000000
$ws.Range("J13").Value = 43000
$ws.Range("K13").Value = 100000000
$ws.Range("L13").Value = 43000
$ws.Range("M13").Value = -99999860
$ws.Range("N13").Value = -43280

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 30007
$ws.Range("J18").Value = 30007
$ws.Range("L18").Value = 30007
$ws.Range("N18").Value = -30353
# Row 23
$ws.Range("H23").Value = 681.9091
$ws.Range("J23").Value = 1000
$ws.Range("L23").Value = 1000
$ws.Range("N23").Value = -1458
# Row 62
$ws.Range("H62").Value = 8495.6
$ws.Range("I62").Value = 8496
$ws.Range("J62").Value = 8495.556
$ws.Range("K62").Value = 8496
$ws.Range("L62").Value = 8495.556
$ws.Range("M62").Value = -7872
$ws.Range("N62").Value = -9743.556
# Row 65
$ws.Range("H65").Value = 8495.6
$ws.Range("I65").Value = 8496
$ws.Range("J65").Value = 8495.556
$ws.Range("K65").Value = 42480
$ws.Range("L65").Value = 42477.78
$ws.Range("M65").Value = -39360
$ws.Range("N65").Value = -48717.78
# Row 88
$ws.Range("H88").Value = 55000
$ws.Range("I88").Value = 55000
$ws.Range("K88").Value = 55000
$ws.Range("M88").Value = -54594
# Row 91
$ws.Range("H91").Value = 55000
$ws.Range("I91").Value = 55000
$ws.Range("K91").Value = 55000
$ws.Range("M91").Value = -53596
# Row 104
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

Write-Host "Applied Seraph_Profits market-price updates"